$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 8.223245666666665
$ws.Range("H2").Value = 24.669737
$ws.Range("I2").Value = 0.003010099901484358
$ws.Range("J2").Value = 0.003010099901484359
$ws.Range("M2").Value = 8.554479333333333
$ws.Range("N2").Value = 25.663438
$ws.Range("O2").Value = 0.1655051910559175
$ws.Range("P2").Value = 0.1655051910559175
$ws.Range("Q2").Value = 70.34558510842288
$ws.Range("R2").Value = 633.1102659758059
$ws.Range("S2").Value = 0.0004981871592925671
$ws.Range("T2").Value = 0.0004981871592925673

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 8.223245666666665
$ws.Range("H3").Value = 24.669737
$ws.Range("I3").Value = 0.003010099901484358
$ws.Range("J3").Value = 0.003010099901484359
$ws.Range("M3").Value = 20.28486166666667
$ws.Range("N3").Value = 60.854585
$ws.Range("O3").Value = 0.392455200938143
$ws.Range("P3").Value = 0.392455200938143
$ws.Range("Q3").Value = 166.8074007993494
$ws.Range("R3").Value = 1501.266607194145
$ws.Range("S3").Value = 0.001181329361680928
$ws.Range("T3").Value = 0.001181329361680928

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 8.223245666666665
$ws.Range("H4").Value = 24.669737
$ws.Range("I4").Value = 0.003010099901484358
$ws.Range("J4").Value = 0.003010099901484359
$ws.Range("M4").Value = 5.037112666666666
$ws.Range("N4").Value = 15.111338
$ws.Range("O4").Value = 0.09745400763531942
$ws.Range("P4").Value = 0.09745400763531943
$ws.Range("Q4").Value = 41.42141490867844
$ws.Range("R4").Value = 372.7927341781059
$ws.Range("S4").Value = 0.0002933462987823308
$ws.Range("T4").Value = 0.000293346298782331

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 8.223245666666665
$ws.Range("H5").Value = 24.669737
$ws.Range("I5").Value = 0.003010099901484358
$ws.Range("J5").Value = 0.003010099901484359
$ws.Range("M5").Value = 17.810622
$ws.Range("N5").Value = 53.431866
$ws.Range("O5").Value = 0.34458560037062
$ws.Range("P5").Value = 0.34458560037062
$ws.Range("Q5").Value = 146.461120182138
$ws.Range("R5").Value = 1318.150081639242
$ws.Range("S5").Value = 0.001037237081728531
$ws.Range("T5").Value = 0.001037237081728532

# Row 6
$ws.Range("I6").Value = 0.9908672065823976
$ws.Range("J6").Value = 0.9908672065823977
$ws.Range("M6").Value = 8.554479333333333
$ws.Range("N6").Value = 25.663438
$ws.Range("O6").Value = 0.1655051910559175
$ws.Range("P6").Value = 0.1655051910559175
$ws.Range("Q6").Value = 23156.41862165932
$ws.Range("R6").Value = 208407.7675949339
$ws.Range("S6").Value = 0.163993666336463
$ws.Range("T6").Value = 0.163993666336463

# Row 7
$ws.Range("I7").Value = 0.9908672065823976
$ws.Range("J7").Value = 0.9908672065823977
$ws.Range("M7").Value = 20.28486166666667
$ws.Range("N7").Value = 60.854585
$ws.Range("O7").Value = 0.392455200938143
$ws.Range("P7").Value = 0.392455200938143
$ws.Range("Q7").Value = 54909.79989927111
$ws.Range("R7").Value = 494188.1990934399
$ws.Range("S7").Value = 0.3888709886623112
$ws.Range("T7").Value = 0.3888709886623113

# Row 8
$ws.Range("I8").Value = 0.9908672065823976
$ws.Range("J8").Value = 0.9908672065823977
$ws.Range("M8").Value = 5.037112666666666
$ws.Range("N8").Value = 15.111338
$ws.Range("O8").Value = 0.09745400763531942
$ws.Range("P8").Value = 0.09745400763531943
$ws.Range("Q8").Value = 13635.13605080458
$ws.Range("R8").Value = 122716.2244572412
$ws.Range("S8").Value = 0.09656398031586859
$ws.Range("T8").Value = 0.09656398031586862

# Row 9
$ws.Range("I9").Value = 0.9908672065823976
$ws.Range("J9").Value = 0.9908672065823977
$ws.Range("M9").Value = 17.810622
$ws.Range("N9").Value = 53.431866
$ws.Range("O9").Value = 0.34458560037062
$ws.Range("P9").Value = 0.34458560037062
$ws.Range("Q9").Value = 48212.19420532844
$ws.Range("R9").Value = 433909.7478479559
$ws.Range("S9").Value = 0.3414385712677546
$ws.Range("T9").Value = 0.3414385712677547

# Row 10
$ws.Range("G10").Value = 14.14340733333333
$ws.Range("H10").Value = 42.430222
$ws.Range("I10").Value = 0.005177161275053701
$ws.Range("J10").Value = 0.005177161275053702
$ws.Range("M10").Value = 8.554479333333333
$ws.Range("N10").Value = 25.663438
$ws.Range("O10").Value = 0.1655051910559175
$ws.Range("P10").Value = 0.1655051910559175
$ws.Range("Q10").Value = 120.9894857359151
$ws.Range("R10").Value = 1088.905371623236
$ws.Range("S10").Value = 0.0008568470659550603
$ws.Range("T10").Value = 0.0008568470659550605

# Row 11
$ws.Range("G11").Value = 14.14340733333333
$ws.Range("H11").Value = 42.430222
$ws.Range("I11").Value = 0.005177161275053701
$ws.Range("J11").Value = 0.005177161275053702
$ws.Range("M11").Value = 20.28486166666667
$ws.Range("N11").Value = 60.854585
$ws.Range("O11").Value = 0.392455200938143
$ws.Range("P11").Value = 0.392455200938143
$ws.Range("Q11").Value = 286.8970612519856
$ws.Range("R11").Value = 2582.07355126787
$ws.Range("S11").Value = 0.002031803868490372
$ws.Range("T11").Value = 0.002031803868490373

# Row 12
$ws.Range("G12").Value = 14.14340733333333
$ws.Range("H12").Value = 42.430222
$ws.Range("I12").Value = 0.005177161275053701
$ws.Range("J12").Value = 0.005177161275053702
$ws.Range("M12").Value = 5.037112666666666
$ws.Range("N12").Value = 15.111338
$ws.Range("O12").Value = 0.09745400763531942
$ws.Range("P12").Value = 0.09745400763531943
$ws.Range("Q12").Value = 71.24193622855955
$ws.Range("R12").Value = 641.177426057036
$ws.Range("S12").Value = 0.0005045351144283633
$ws.Range("T12").Value = 0.0005045351144283636

# Row 13
$ws.Range("G13").Value = 14.14340733333333
$ws.Range("H13").Value = 42.430222
$ws.Range("I13").Value = 0.005177161275053701
$ws.Range("J13").Value = 0.005177161275053702
$ws.Range("M13").Value = 17.810622
$ws.Range("N13").Value = 53.431866
$ws.Range("O13").Value = 0.34458560037062
$ws.Range("P13").Value = 0.34458560037062
$ws.Range("Q13").Value = 251.902881806028
$ws.Range("R13").Value = 2267.125936254252
$ws.Range("S13").Value = 0.001783975226179904
$ws.Range("T13").Value = 0.001783975226179904

# Row 14
$ws.Range("G14").Value = 2.583085
$ws.Range("H14").Value = 7.749255
$ws.Range("I14").Value = 0.0009455322410643118
$ws.Range("J14").Value = 0.0009455322410643119
$ws.Range("M14").Value = 8.554479333333333
$ws.Range("N14").Value = 25.663438
$ws.Range("O14").Value = 0.1655051910559175
$ws.Range("P14").Value = 0.1655051910559175
$ws.Range("Q14").Value = 22.09694724874333
$ws.Range("R14").Value = 198.87252523869
$ws.Range("S14").Value = 0.0001564904942068788
$ws.Range("T14").Value = 0.0001564904942068788

# Row 15
$ws.Range("G15").Value = 2.583085
$ws.Range("H15").Value = 7.749255
$ws.Range("I15").Value = 0.0009455322410643118
$ws.Range("J15").Value = 0.0009455322410643119
$ws.Range("M15").Value = 20.28486166666667
$ws.Range("N15").Value = 60.854585
$ws.Range("O15").Value = 0.392455200938143
$ws.Range("P15").Value = 0.392455200938143
$ws.Range("Q15").Value = 52.39752189824167
$ws.Range("R15").Value = 471.577697084175
$ws.Range("S15").Value = 0.0003710790456603871
$ws.Range("T15").Value = 0.0003710790456603872

# Row 16
$ws.Range("G16").Value = 2.583085
$ws.Range("H16").Value = 7.749255
$ws.Range("I16").Value = 0.0009455322410643118
$ws.Range("J16").Value = 0.0009455322410643119
$ws.Range("M16").Value = 5.037112666666666
$ws.Range("N16").Value = 15.111338
$ws.Range("O16").Value = 0.09745400763531942
$ws.Range("P16").Value = 0.09745400763531943
$ws.Range("Q16").Value = 13.01129017257667
$ws.Range("R16").Value = 117.10161155319
$ws.Range("S16").Value = 0.00009214590624012212
$ws.Range("T16").Value = 0.00009214590624012214

# Row 17
$ws.Range("G17").Value = 2.583085
$ws.Range("H17").Value = 7.749255
$ws.Range("I17").Value = 0.0009455322410643118
$ws.Range("J17").Value = 0.0009455322410643119
$ws.Range("M17").Value = 17.810622
$ws.Range("N17").Value = 53.431866
$ws.Range("O17").Value = 0.34458560037062
$ws.Range("P17").Value = 0.34458560037062
$ws.Range("Q17").Value = 46.00635052886999
$ws.Range("R17").Value = 414.05715475983
$ws.Range("S17").Value = 0.0003258167949569236
$ws.Range("T17").Value = 0.0003258167949569237
